# Update the dSF column (F) values for specific rows as part of a data repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = -7
$ws.Range("F10").Value = -6
$ws.Range("F13").Value = -2
$ws.Range("F14").Value = -5
$ws.Range("F17").Value = 2
